$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.730.32"
$ws.Range("E2").Value = "  +5.32%  "
$ws.Range("D3").Value = "'2.302.49"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").Value = "'302.06"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "'101.58"
$ws.Range("E6").Value = "  +12.79%  "
$ws.Range("D7").Value = "'0.570"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  +8.53%  "
$ws.Range("D10").Value = "'36.83"
$ws.Range("E10").Value = "  +11.25%  "
$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").Value = "'7.37"
$ws.Range("E12").Value = "  +5.79%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'2.653.99"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").Value = "'2.308.29"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("E17").Value = "  +4.47%  "
$ws.Range("D18").Value = "'46.705.53"
$ws.Range("E18").Value = "  +5.87%  "
$ws.Range("D19").Value = "'13.49"
$ws.Range("E19").Value = "  +21.52%  "
$ws.Range("D20").Value = "'0.0₃0946"
$ws.Range("E20").Value = "  +4.30%  "
$ws.Range("D21").Value = "'6.11"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").Value = "'66.76"
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").Value = "'248.52"
$ws.Range("E23").Value = "  +5.41%  "
$ws.Range("E24").Value = "  +5.54%  "
$ws.Range("E25").Value = "  +5.57%  "
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("D27").Value = "'44.95"
$ws.Range("E27").Value = "  +16.40%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Value = "'9.90"
$ws.Range("E29").Value = "  +5.84%  "
$ws.Range("D30").Value = "'20.14"
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("D31").Value = "'5.79"
$ws.Range("E31").Value = "  +7.44%  "
$ws.Range("D32").Value = "'147.41"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "'0.0800"
$ws.Range("E33").Value = "  +6.75%  "
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("E35").Value = "  +10.54%  "
$ws.Range("E36").Value = "  +9.46%  "
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("E38").Value = "  +8.58%  "
$ws.Range("E39").Value = "  +20.65%  "
$ws.Range("E40").Value = "  +13.83%  "
$ws.Range("D41").Value = "'3.52"
$ws.Range("E41").Value = "  +10.31%  "
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'1.867.80"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").Value = "'1.97"
$ws.Range("E45").Value = "  +11.37%  "
$ws.Range("D46").Value = "'87.87"
$ws.Range("E46").Value = "  +19.09%  "
$ws.Range("D47").Value = "'0.196"
$ws.Range("E47").Value = "  +9.74%  "
$ws.Range("D48").Value = "'74.31"
$ws.Range("E48").Value = "  +10.07%  "
$ws.Range("D49").Value = "'4.91"
$ws.Range("E49").Value = "  +10.77%  "
$ws.Range("D50").Value = "'97.30"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("D51").Value = "'8.08"
$ws.Range("E51").Value = "  +5.65%  "
